# Chatterbox "text to speech" edit
#
# 1) The editing-session bookmark "_GoBack" (Word drops this at the spot of
#    the last edit) moves from the Title paragraph ("Chatterbox") down to the
#    last bullet of Iteration 3 ("AI can hold a conversation with you") -
#    i.e. that is where the author's final edit of this session landed.
# 2) Word silently minted the built-in "Balloon Text" / "Balloon Text Char"
#    style pair into the style sheet (a side-effect of the editing session -
#    e.g. turning on a reviewing/"read aloud" style feature while adding the
#    text-to-speech bullet).

$d = $word.ActiveDocument

# --- 1. Relocate the "_GoBack" bookmark ------------------------------------

# Remove it from wherever it currently lives (the Title paragraph).
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# Find the new insertion point - the start of the "AI can hold a
# conversation with you" bullet - and drop a fresh, collapsed "_GoBack"
# bookmark right before its text run.
$target = $d.Content
$found = $target.Find.Execute("AI can hold a conversation with you", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $target.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $target)
}

# --- 2. Add the "Balloon Text" / "Balloon Text Char" style pair ------------

$balloonText = $d.Styles.Add("Balloon Text", 1)
$balloonText.BaseStyle = "Normal"
$balloonText.LinkStyle = "BalloonTextChar"
$balloonText.Priority = 99
$balloonText.UnhideWhenUsed = $true
$balloonText.Font.Name = "Tahoma"
$balloonText.Font.NameBi = "Tahoma"
$balloonText.Font.Size = 8
$balloonText.Font.SizeBi = 8
$balloonText.ParagraphFormat.SpaceAfter = 0
$balloonText.ParagraphFormat.LineSpacingRule = 0
$balloonText.ParagraphFormat.LineSpacing = 12

$balloonTextChar = $d.Styles.Add("Balloon Text Char", 2)
$balloonTextChar.BaseStyle = "DefaultParagraphFont"
$balloonTextChar.LinkStyle = "BalloonText"
$balloonTextChar.Priority = 99
$balloonTextChar.Font.Name = "Tahoma"
$balloonTextChar.Font.NameBi = "Tahoma"
$balloonTextChar.Font.Size = 8
$balloonTextChar.Font.SizeBi = 8
